$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 22: task 19 - mark done, progress 100%
$ws.Cells.Item(22,5).Value = 1
$ws.Range("G7").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Cells.Item(22,7).Value = "done"

# Row 26: task 23 - mark done, progress 100%
$ws.Cells.Item(26,5).Value = 1
$ws.Range("G7").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Cells.Item(26,7).Value = "done"

# Row 27: task 24 - mark cancel
$ws.Range("G4").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Cells.Item(27,7).Value = "cancel"

$excel.CutCopyMode = $false
